$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
    @("Nikola Jovic", "PF,C", "Miami Heat"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Mike Conley", "PG", "Minnesota Timberwolves"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
    @("Jerami Grant", "SF,PF", "Portland Trail Blazers"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jimmy Butler", "SF,PF", "Miami Heat")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
